# Applies the commit "add extra regions and models, turn off defective check"
# to available_models_regions_variables_units.xlsx
#
#  * variable_units sheet: 15 new Variable/Unit rows (1815-1829) for
#    GEMINI-E3 related variables that didn't have a unit entry yet.
#  * models sheet: 2 new model rows (51-52) for GEMINI-E3 8.0 / 7.0.
#  * regions sheet: 25 new region rows (266-290) used by GEMINI-E3.
#  * widen the "duplicate value" conditional-format range on variable_units
#    to cover the freshly added rows (plus a couple of spare rows, matching
#    the author's manual drag-fill).

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("variable_units")
$wsModels    = $wb.Worksheets.Item("models")
$wsRegions   = $wb.Worksheets.Item("regions")

# ---------------------------------------------------------------------
# models: two new GEMINI-E3 model versions
# ---------------------------------------------------------------------
$newModels = @("GEMINI-E3 8.0", "GEMINI-E3 7.0")
$row = 51
foreach ($m in $newModels) {
    $wsModels.Cells.Item($row, 1).Value = $m
    $row++
}

# ---------------------------------------------------------------------
# regions: GEMINI-E3 region list
# ---------------------------------------------------------------------
$newRegions = @(
    "DEU", "FRA", "ITA", "SPN", "NLD", "SWE", "POL", "BEL",
    "EU1", "EU2", "EU3", "EU4", "EU5", "EU6",
    "GBR", "USA", "CHI", "IND", "RUS", "CSA", "MID", "AFR", "ASI", "ROW", "WORLD"
)
$row = 266
foreach ($rg in $newRegions) {
    $wsRegions.Cells.Item($row, 1).Value = $rg
    $row++
}

# ---------------------------------------------------------------------
# variable_units: new Variable / Unit rows
# ---------------------------------------------------------------------
$unitRows = @(
    @{ Row = 1815; Variable = "Capital Stock";                             Unit = "billion US$2010/yr or local currency/yr "; Black = $false },
    @{ Row = 1816; Variable = "Expenditure|government";                    Unit = "billion US$2010/yr OR local currency";      Black = $false },
    @{ Row = 1817; Variable = "Expenditure|household";                     Unit = "billion US$2010/yr OR local currency";      Black = $false },
    @{ Row = 1818; Variable = "Expenditure|household|Energy";               Unit = "billion US$2010/yr OR local currency";      Black = $false },
    @{ Row = 1819; Variable = "Export";                                    Unit = "billion US$2010/yr OR local currency";      Black = $false },
    @{ Row = 1820; Variable = "GDP|MER";                                   Unit = "billion US$2010/yr OR local currency";      Black = $false },
    @{ Row = 1821; Variable = "GDP|PPP";                                   Unit = "billion US$2010/yr OR local currency";      Black = $true  },
    @{ Row = 1822; Variable = "Import";                                    Unit = "billion US$2010/yr OR local currency";      Black = $true  },
    @{ Row = 1823; Variable = "Investment";                                Unit = "billion US$2010/yr OR local currency";      Black = $true  },
    @{ Row = 1824; Variable = "Investment|Energy Supply";                  Unit = "billion US$2010/yr OR local currency";      Black = $true  },
    @{ Row = 1825; Variable = "Investment|Energy Supply|Electricity";      Unit = "billion US$2010/yr or local currency/yr";   Black = $false },
    @{ Row = 1826; Variable = "Price|Primary Energy|Oil";                  Unit = "US$2010/GJ or local currency/GJ";           Black = $false },
    @{ Row = 1827; Variable = "Value Added|Agriculture";                   Unit = "billion US$2010/yr OR local currency";      Black = $true  },
    @{ Row = 1828; Variable = "Value Added|Industry|Energy";               Unit = "billion US$2010/yr OR local currency";      Black = $true  },
    @{ Row = 1829; Variable = "Value Added|Industry|Energy Intensive";     Unit = "billion US$2010/yr OR local currency";      Black = $true  }
)

# Column A first (every label here already exists in the shared-string
# table from elsewhere on the sheet, so the write order has no effect on
# the resulting shared-string table).
foreach ($r in $unitRows) {
    $wsVariables.Cells.Item($r.Row, 1).Value = $r.Variable
}

# Column B next, in the order the four distinct unit strings were first
# introduced (row 1826's unit text was entered before row 1825's, so the
# shared-string table picks it up first).
$unitColumnOrder = @(1815, 1816, 1817, 1818, 1819, 1820, 1821, 1822, 1823, 1824, 1826, 1825, 1827, 1828, 1829)
foreach ($rowNum in $unitColumnOrder) {
    $r = $unitRows | Where-Object { $_.Row -eq $rowNum }
    $cellB = $wsVariables.Cells.Item($r.Row, 2)
    $cellB.Value = $r.Unit
    if ($r.Black) {
        # A handful of the pasted-in rows carry an explicit black font
        # color (as opposed to the sheet's default automatic/theme color).
        $cellB.Font.Color = 0
    }
}

# Widen the duplicate-value conditional formatting so it keeps covering
# column A after the new rows were appended (the author's manual fill
# handle drag went a couple of rows past the last data row).
$fcs = $wsVariables.Range("A2:A1814").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($wsVariables.Range("A2:A1831"))

# ---------------------------------------------------------------------
# Mirror the author's final cell selections on each sheet, then leave
# variable_units as the active/visible tab (as in the saved file).
# ---------------------------------------------------------------------
$wsModels.Activate()
$wsModels.Range("G36").Select()

$wsRegions.Activate()
$wsRegions.Range("D273").Select()

$wsVariables.Activate()
$wsVariables.Range("A1830:B1835").Select()
